$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# COSTR row: remove check (more datapoints/BU/experimentos counted)
$ws.Range("C2").Value = "'82132"
$ws.Range("D2").Value = "'343"
$ws.Range("E2").Value = "'701"

# GRLSR row
$ws.Range("C3").Value = "'11149"
$ws.Range("D3").Value = "'69"
$ws.Range("E3").Value = "'124"

# HELMR row
$ws.Range("C4").Value = "'26374"
$ws.Range("D4").Value = "'157"
$ws.Range("E4").Value = "'351"

# LFSPR row
$ws.Range("C5").Value = "'36752"
$ws.Range("D5").Value = "'205"
$ws.Range("E5").Value = "'376"

# SCLBR row
$ws.Range("C6").Value = "'12620"
$ws.Range("D6").Value = "'64"
$ws.Range("E6").Value = "'121"

# DLLFR row
$ws.Range("C8").Value = "'7060"
$ws.Range("D8").Value = "'80"
$ws.Range("E8").Value = "'133"

# PRMDN row
$ws.Range("C9").Value = "'8914"
$ws.Range("D9").Value = "'42"
$ws.Range("E9").Value = "'76"

# Total row
$ws.Range("C10").Value = "'187640"
$ws.Range("D10").Value = "'995"
$ws.Range("E10").Value = "'1920"
